$wb = $excel.ActiveWorkbook

# ALC (sheet1) row 2
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 439.75
$ws.Range("I2").Value = 350
$ws.Range("J2").Value = 529.5
$ws.Range("K2").Value = 350
$ws.Range("L2").Value = 529.5
$ws.Range("M2").Value = -237
$ws.Range("N2").Value = -755.5

# ALC (sheet1) row 3
$ws = $wb.Worksheets.Item(1)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# ALC (sheet1) row 102
$ws = $wb.Worksheets.Item(1)
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# ALC (sheet1) row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 1310.4048
$ws.Range("I132").Value = 1238.425
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 3715.275
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -1185.275
$ws.Range("N132").Value = -13310

# ARM (sheet2) row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 40370.074
$ws.Range("I32").Value = 43495.72
$ws.Range("J32").Value = 1299.5
$ws.Range("K32").Value = 43495.72
$ws.Range("L32").Value = 1299.5
$ws.Range("M32").Value = -43208.72
$ws.Range("N32").Value = -1873.5

# ARM (sheet2) row 41
$ws = $wb.Worksheets.Item(2)
$ws.Range("H41").Value = 8750
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = -4586
$ws.Range("N41").Value = -10828

# ARM (sheet2) row 61
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 10998.875
$ws.Range("I61").Value = 6960.9414
$ws.Range("J61").Value = 15575.2
$ws.Range("K61").Value = 6960.9414
$ws.Range("L61").Value = 15575.2
$ws.Range("M61").Value = -6748.9414
$ws.Range("N61").Value = -15999.2

# ARM (sheet2) row 102
$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 2488.889
$ws.Range("I102").Value = 2116.6667
$ws.Range("J102").Value = 3233.3333
$ws.Range("K102").Value = 2116.6667
$ws.Range("L102").Value = 3233.3333
$ws.Range("M102").Value = -494.6667000000002
$ws.Range("N102").Value = -6477.3333

# ARM (sheet2) row 113
$ws = $wb.Worksheets.Item(2)
$ws.Range("H113").Value = 76399.75
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 76399.75
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 76399.75
$ws.Range("N113").Value = -85077.75

# ARM (sheet2) row 122
$ws = $wb.Worksheets.Item(2)
$ws.Range("H122").Value = 1875.1428
$ws.Range("I122").Value = 1729.3334
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 5188.0002
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -2738.0002
$ws.Range("N122").Value = -13150

# ARM (sheet2) row 132
$ws = $wb.Worksheets.Item(2)
$ws.Range("H132").Value = 9149.200000000001
$ws.Range("I132").Value = 10017.23
$ws.Range("J132").Value = 3507
$ws.Range("K132").Value = 30051.69
$ws.Range("L132").Value = 10521
$ws.Range("M132").Value = -27521.69
$ws.Range("N132").Value = -15581

# ARM (sheet2) row 135
$ws = $wb.Worksheets.Item(2)
$ws.Range("H135").Value = 42845.8
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 42845.8
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 42845.8
$ws.Range("N135").Value = -52985.8

# ARM (sheet2) row 136
$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 10998.875
$ws.Range("I136").Value = 6960.9414
$ws.Range("J136").Value = 15575.2
$ws.Range("K136").Value = 20882.8242
$ws.Range("L136").Value = 46725.60000000001
$ws.Range("M136").Value = -18332.8242
$ws.Range("N136").Value = -51825.60000000001

# BSM (sheet3) row 80
$ws = $wb.Worksheets.Item(3)
$ws.Range("H80").Value = 179.72728
$ws.Range("I80").Value = 115
$ws.Range("J80").Value = 204
$ws.Range("K80").Value = 115
$ws.Range("L80").Value = 204
$ws.Range("M80").Value = 883
$ws.Range("N80").Value = -2200

# BSM (sheet3) row 83
$ws = $wb.Worksheets.Item(3)
$ws.Range("H83").Value = 179.72728
$ws.Range("I83").Value = 115
$ws.Range("J83").Value = 204
$ws.Range("K83").Value = 575
$ws.Range("L83").Value = 1020
$ws.Range("M83").Value = 4417
$ws.Range("N83").Value = -11004

# BSM (sheet3) row 86
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 1758.7255
$ws.Range("I86").Value = 1744.0625
$ws.Range("J86").Value = 1993.3334
$ws.Range("K86").Value = 1744.0625
$ws.Range("L86").Value = 1993.3334
$ws.Range("M86").Value = -621.0625
$ws.Range("N86").Value = -4239.3334

# BSM (sheet3) row 89
$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value = 1758.7255
$ws.Range("I89").Value = 1744.0625
$ws.Range("J89").Value = 1993.3334
$ws.Range("K89").Value = 8720.3125
$ws.Range("L89").Value = 9966.666999999999
$ws.Range("M89").Value = -3104.3125
$ws.Range("N89").Value = -21198.667

# CRP (sheet4) row 82
$ws = $wb.Worksheets.Item(4)
$ws.Range("H82").Value = 35480.668
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 35480.668
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 35480.668
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -36202.668

# CRP (sheet4) row 85
$ws = $wb.Worksheets.Item(4)
$ws.Range("H85").Value = 35480.668
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 35480.668
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 35480.668
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -37976.668

# CRP (sheet4) row 86
$ws = $wb.Worksheets.Item(4)
$ws.Range("H86").Value = 3534.5
$ws.Range("I86").Value = 3541.4
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 3541.4
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -2418.4
$ws.Range("N86").Value = -5746

# CRP (sheet4) row 89
$ws = $wb.Worksheets.Item(4)
$ws.Range("H89").Value = 3534.5
$ws.Range("I89").Value = 3541.4
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 17707
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -12091
$ws.Range("N89").Value = -28732

# CUL (sheet5) row 46
$ws = $wb.Worksheets.Item(5)
$ws.Range("H46").Value = 2835.6323
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2863.0298
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 8589.089399999999
$ws.Range("M46").Value = -2909
$ws.Range("N46").Value = -8771.089399999999

# CUL (sheet5) row 102
$ws = $wb.Worksheets.Item(5)
$ws.Range("H102").Value = 4931.9
$ws.Range("I102").Value = 3013
$ws.Range("J102").Value = 5411.625
$ws.Range("K102").Value = 9039
$ws.Range("L102").Value = 16234.875
$ws.Range("M102").Value = -6605
$ws.Range("N102").Value = -21102.875

# CUL (sheet5) row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 835.5
$ws.Range("I131").Value = 661.125
$ws.Range("J131").Value = 975
$ws.Range("K131").Value = 1983.375
$ws.Range("L131").Value = 2925
$ws.Range("M131").Value = 3056.625
$ws.Range("N131").Value = -13005

# CUL (sheet5) row 132
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 1709.1428
$ws.Range("I132").Value = 1720.7273
$ws.Range("J132").Value = 1666.6666
$ws.Range("K132").Value = 15486.5457
$ws.Range("L132").Value = 14999.9994
$ws.Range("M132").Value = -12956.5457
$ws.Range("N132").Value = -20059.9994

# GSM (sheet6) row 41
$ws = $wb.Worksheets.Item(6)
$ws.Range("H41").Value = 7183.6665
$ws.Range("I41").Value = 2275.5
$ws.Range("J41").Value = 17000
$ws.Range("K41").Value = 2275.5
$ws.Range("L41").Value = 17000
$ws.Range("M41").Value = -1920.5
$ws.Range("N41").Value = -17710

# GSM (sheet6) row 70
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 5559.3257
$ws.Range("I70").Value = 5325.8965
$ws.Range("J70").Value = 6042.857
$ws.Range("K70").Value = 5325.8965
$ws.Range("L70").Value = 6042.857
$ws.Range("M70").Value = -5055.8965
$ws.Range("N70").Value = -6582.857

# GSM (sheet6) row 73
$ws = $wb.Worksheets.Item(6)
$ws.Range("H73").Value = 5559.3257
$ws.Range("I73").Value = 5325.8965
$ws.Range("J73").Value = 6042.857
$ws.Range("K73").Value = 5325.8965
$ws.Range("L73").Value = 6042.857
$ws.Range("M73").Value = -4389.8965
$ws.Range("N73").Value = -7914.857

# GSM (sheet6) row 126
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 3012.5881
$ws.Range("I126").Value = 1977.7778
$ws.Range("J126").Value = 4176.75
$ws.Range("K126").Value = 5933.3334
$ws.Range("L126").Value = 12530.25
$ws.Range("M126").Value = -3463.3334
$ws.Range("N126").Value = -17470.25

# LTW (sheet7) row 100
$ws = $wb.Worksheets.Item(7)
$ws.Range("H100").Value = 4212.5
$ws.Range("I100").Value = 4100
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 4100
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -3559
$ws.Range("N100").Value = -6082

# LTW (sheet7) row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 7871.645
$ws.Range("I122").Value = 7425.8
$ws.Range("J122").Value = 8682.272000000001
$ws.Range("K122").Value = 22277.4
$ws.Range("L122").Value = 26046.816
$ws.Range("M122").Value = -19827.4
$ws.Range("N122").Value = -30946.816

# LTW (sheet7) row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 5959.1934
$ws.Range("I136").Value = 3310.3845
$ws.Range("J136").Value = 7872.222
$ws.Range("K136").Value = 9931.1535
$ws.Range("L136").Value = 23616.666
$ws.Range("M136").Value = -7381.1535
$ws.Range("N136").Value = -28716.666

# WVR (sheet8) row 47
$ws = $wb.Worksheets.Item(8)
$ws.Range("H47").Value = 300000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 300000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 300000
$ws.Range("N47").Value = -301144

# WVR (sheet8) row 122
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 5198.3213
$ws.Range("I122").Value = 1597.8
$ws.Range("J122").Value = 14199.625
$ws.Range("K122").Value = 4793.4
$ws.Range("L122").Value = 42598.875
$ws.Range("M122").Value = -2343.4
$ws.Range("N122").Value = -47498.875

# WVR (sheet8) row 126
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 1800.875
$ws.Range("I126").Value = 1846.4615
$ws.Range("J126").Value = 1603.3334
$ws.Range("K126").Value = 5539.3845
$ws.Range("L126").Value = 4810.0002
$ws.Range("M126").Value = -3069.3845
$ws.Range("N126").Value = -9750.0002

# WVR (sheet8) row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 4959.0786
$ws.Range("I136").Value = 2148.818
$ws.Range("J136").Value = 10111.223
$ws.Range("K136").Value = 6446.454000000001
$ws.Range("L136").Value = 30333.669
$ws.Range("M136").Value = -3896.454000000001
$ws.Range("N136").Value = -35433.669
